$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume 1h (E) columns with refreshed crypto data.
# D-column values are forced to text (matching the original inlineStr cells)
# by temporarily applying a text NumberFormat before assignment, then clearing
# the format again so the cell keeps its original (unstyled) appearance.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.587.47'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.444.06'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.95%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '592.91'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.64%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.60'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.67%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.444.21'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -3.00%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.492'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.39'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -6.09%  '
$ws.Range('E11').Value = '  -8.17%  '
$ws.Range('E12').Value = '  -7.29%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.019.64'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.22%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000181'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -10.07%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '26.67'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -8.35%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.430.15'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.49%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '65.482.61'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('E18').Value = '  -2.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.97'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -8.99%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.82'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.53%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.70'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -6.63%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '394.41'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.66%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.548'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -9.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '73.40'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -5.52%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.581.30'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.08%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000106'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -9.68%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.24'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -8.03%  '
$ws.Range('E30').Value = '  -8.42%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.21'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -10.38%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.444.58'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.92%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  -6.13%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '23.04'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.37%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '171.61'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.74%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.95'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -8.89%  '
$ws.Range('E38').Value = '  -12.04%  '
$ws.Range('E39').Value = '  -6.52%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.83'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -9.84%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0770'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -7.10%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.825'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.89%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '43.57'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.77%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.999'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.44'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -12.98%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.63'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -10.64%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.11'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '22.74'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.57'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.91%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.11'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -14.07%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.205.16'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.80%  '
